$d = $word.ActiveDocument

$newParagraphs = @(
    "StudentDtoRequest(id=17, birthday=2007-05-30, parentStatus=FATHER, userId=1, parentId=2, gradeId=2)",
    "1 StudentDtoRequest {id=17, birthday=2007-05-30, parentStatus=FATHER, userId=1, parentId=2, gradeId=2}",
    "1 StudentDtoRequest {id=2, birthday=2007-05-30, parentStatus=FATHER, userId=1, parentId=2, gradeId=2}"
)

foreach ($text in $newParagraphs) {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $r.Collapse(0)
    $r.InsertAfter($text)
}
